$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 02:03"

# Estados Unidos (row 4) - updated totals
$ws.Range("B4").Value = 1237045
$ws.Range("C4").Value = 24210
$ws.Range("D4").Value = 199691
$ws.Range("E4").Value = 965112
$ws.Range("G4").Value = 2321
$ws.Range("H4").Value = 72242

# Argentina (row 57) - updated totals
$ws.Range("B57").Value = 5020
$ws.Range("C57").Value = 133
$ws.Range("E57").Value = 3284
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 264

# Curazao overtakes Dominica (rows 198-199 swap order/values)
$ws.Range("A198").Value = "Curazao"
$ws.Range("D198").Value = 13
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "Dominica"
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 0

# Seychelles overtakes Montserrat (rows 205-206 swap order/values)
$ws.Range("A205").Value = "Seychelles"
$ws.Range("D205").Value = 8
$ws.Range("F205").Value = 0
$ws.Range("H205").Value = 0

$ws.Range("A206").Value = "Montserrat"
$ws.Range("D206").Value = 7
$ws.Range("F206").Value = 1
$ws.Range("H206").Value = 1
